# followed_users.xlsx edit: new rows of followed users + selection/window bookkeeping

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove workbook protection element (<workbookProtection/> dropped in target)
$wb.Unprotect()

# New table contents: 12 rows x 2 columns (username, timestamp) replacing the
# old 4-row sample data.
$rows = @(
    @("alash_arts", "01/09/2020, 11:59:37"),
    @("lolariostyle", "01/09/2020, 11:59:50"),
    @("yleniastorti", "01/09/2020, 12:06:04"),
    @("swissmonamour", "01/09/2020, 12:06:11"),
    @("fashiongirls91", "01/09/2020, 12:06:30"),
    @("mesiszigeti", "01/09/2020, 12:06:37"),
    @("lolariostyle", "01/09/2020, 12:22:55"),
    @("fashiongirls91", "01/09/2020, 12:23:01"),
    @("mesiszigeti", "01/09/2020, 12:23:08"),
    @("ikizlerlekurabiyeler", "01/09/2020, 12:31:29"),
    @("_.fayis2", "01/09/2020, 12:31:36"),
    @("italia_dev", "01/09/2020, 12:31:42")
)

$rowCount = $rows.Count
$arr = New-Object 'object[,]' $rowCount,2
for ($i = 0; $i -lt $rowCount; $i++) {
    $arr[$i,0] = $rows[$i][0]
    $arr[$i,1] = $rows[$i][1]
}

$ws.Range("A1:B12").Value = $arr

# Update the sheet selection to C10 and make sure the sheet is marked active/selected
$ws.Range("C10").Select() | Out-Null
